$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.378.66"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.722.44"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.23"
$ws.Range("E5").Value = "  -1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4869"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2589"
$ws.Range("E8").Value = "  -3.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06185"
$ws.Range("D10").Value = "1.722.95"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06975"
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.50"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.518"
$ws.Range("E13").Value = "  -2.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5973"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.11"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9995"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "26.384.31"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9996"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007184"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.30"
$ws.Range("E20").Value = "  -1.73%  "
$ws.Range("D21").Value = "1.945.28"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.429"
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.485"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.089"
$ws.Range("E24").Value = "  -3.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.85"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.22"
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.29"
$ws.Range("E28").Value = "  -1.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.723"
$ws.Range("E29").Value = "  -3.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.898"
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08003"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.647"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04492"
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9954"
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6213"
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9283"
$ws.Range("E37").Value = "  +2.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.955"
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.385"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9991"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01476"
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.55"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.433"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3834"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.867"
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.18"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.650"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "50.98"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.217"
$ws.Range("E51").Value = "  -2.48%  "
